# NID ajuste ... 01122023
# Swap the presenter-name order on the title slide:
#   "Palestrante(s): Josiane Sheila e Júlia Cardoso"
#   -> "Palestrante(s): Júlia Cardoso e Josiane Sheila"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the textbox that holds the "Palestrante(s): ..." line instead of
# hard-coding a shape index.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
        if ($candidate.TextFrame.TextRange.Text -like "Palestrante(s)*") {
            $shp = $candidate
        }
    }
}

$tr = $shp.TextFrame.TextRange

# Only retype the part of the text that actually changes ("): <names>"),
# leaving the leading "Palestrante(s" run untouched - this mirrors how the
# author edited just the trailing portion of the line in PowerPoint,
# which splits the paragraph into two runs.
$splitPoint = 14
$tail = $tr.Characters($splitPoint, $tr.Length - ($splitPoint - 1))
$tail.Text = "): Júlia Cardoso e Josiane Sheila"
